$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.218698
$ws.Range("H2").Value = 30.656094
$ws.Range("I2").Value = 0.01131191978527373
$ws.Range("J2").Value = 0.01131191978527373
$ws.Range("M2").Value = 1.918906333333333
$ws.Range("N2").Value = 5.756718999999999
$ws.Range("O2").Value = 0.006524019162508824
$ws.Range("P2").Value = 0.006524019162508824
$ws.Range("Q2").Value = 19.60872431062066
$ws.Range("R2").Value = 176.478518795586
$ws.Range("S2").Value = 0.00007379918144388852
$ws.Range("T2").Value = 0.00007379918144388854
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.218698
$ws.Range("H3").Value = 30.656094
$ws.Range("I3").Value = 0.01131191978527373
$ws.Range("J3").Value = 0.01131191978527373
$ws.Range("O3").Value = 0.6163557430885885
$ws.Range("P3").Value = 0.6163557430885885
$ws.Range("Q3").Value = 1852.531321941151
$ws.Range("R3").Value = 16672.78189747036
$ws.Range("S3").Value = 0.006972166725010898
$ws.Range("T3").Value = 0.006972166725010899
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.218698
$ws.Range("H4").Value = 30.656094
$ws.Range("I4").Value = 0.01131191978527373
$ws.Range("J4").Value = 0.01131191978527373
$ws.Range("M4").Value = 29.04767233333333
$ws.Range("N4").Value = 87.143017
$ws.Range("O4").Value = 0.09875811426384234
$ws.Range("P4").Value = 0.09875811426384236
$ws.Range("Q4").Value = 296.8293911772886
$ws.Range("R4").Value = 2671.464520595598
$ws.Range("S4").Value = 0.001117143866697482
$ws.Range("T4").Value = 0.001117143866697483
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.218698
$ws.Range("H5").Value = 30.656094
$ws.Range("I5").Value = 0.01131191978527373
$ws.Range("J5").Value = 0.01131191978527373
$ws.Range("M5").Value = 81.87450533333333
$ws.Range("N5").Value = 245.623516
$ws.Range("O5").Value = 0.2783621234850603
$ws.Range("P5").Value = 0.2783621234850603
$ws.Range("Q5").Value = 836.6508439007225
$ws.Range("R5").Value = 7529.857595106503
$ws.Range("S5").Value = 0.003148810012121463
$ws.Range("T5").Value = 0.003148810012121464
$ws.Range("I6").Value = 0.8540056659644313
$ws.Range("J6").Value = 0.8540056659644314
$ws.Range("M6").Value = 1.918906333333333
$ws.Range("N6").Value = 5.756718999999999
$ws.Range("O6").Value = 0.006524019162508824
$ws.Range("P6").Value = 0.006524019162508824
$ws.Range("Q6").Value = 1480.381931757069
$ws.Range("R6").Value = 13323.43738581362
$ws.Range("S6").Value = 0.00557154932964306
$ws.Range("T6").Value = 0.005571549329643061
$ws.Range("I7").Value = 0.8540056659644313
$ws.Range("J7").Value = 0.8540056659644314
$ws.Range("O7").Value = 0.6163557430885885
$ws.Range("P7").Value = 0.6163557430885885
$ws.Range("S7").Value = 0.5263712968473719
$ws.Range("T7").Value = 0.526371296847372
$ws.Range("I8").Value = 0.8540056659644313
$ws.Range("J8").Value = 0.8540056659644314
$ws.Range("M8").Value = 29.04767233333333
$ws.Range("N8").Value = 87.143017
$ws.Range("O8").Value = 0.09875811426384234
$ws.Range("P8").Value = 0.09875811426384236
$ws.Range("Q8").Value = 22409.45716572219
$ws.Range("R8").Value = 201685.1144914997
$ws.Range("S8").Value = 0.08433998914128409
$ws.Range("T8").Value = 0.0843399891412841
$ws.Range("I9").Value = 0.8540056659644313
$ws.Range("J9").Value = 0.8540056659644314
$ws.Range("M9").Value = 81.87450533333333
$ws.Range("N9").Value = 245.623516
$ws.Range("O9").Value = 0.2783621234850603
$ws.Range("P9").Value = 0.2783621234850603
$ws.Range("Q9").Value = 63163.86384345723
$ws.Range("R9").Value = 568474.7745911151
$ws.Range("S9").Value = 0.2377228306461322
$ws.Range("T9").Value = 0.2377228306461322
$ws.Range("G10").Value = 121.2114333333333
$ws.Range("H10").Value = 363.6343
$ws.Range("I10").Value = 0.1341789346279459
$ws.Range("J10").Value = 0.1341789346279459
$ws.Range("M10").Value = 1.918906333333333
$ws.Range("N10").Value = 5.756718999999999
$ws.Range("O10").Value = 0.006524019162508824
$ws.Range("P10").Value = 0.006524019162508824
$ws.Range("Q10").Value = 232.5933870957444
$ws.Range("R10").Value = 2093.3404838617
$ws.Range("S10").Value = 0.0008753859407177378
$ws.Range("T10").Value = 0.0008753859407177378
$ws.Range("G11").Value = 121.2114333333333
$ws.Range("H11").Value = 363.6343
$ws.Range("I11").Value = 0.1341789346279459
$ws.Range("J11").Value = 0.1341789346279459
$ws.Range("O11").Value = 0.6163557430885885
$ws.Range("P11").Value = 0.6163557430885885
$ws.Range("Q11").Value = 21974.22576020758
$ws.Range("R11").Value = 197768.0318418682
$ws.Range("S11").Value = 0.08270195695944273
$ws.Range("T11").Value = 0.08270195695944273
$ws.Range("G12").Value = 121.2114333333333
$ws.Range("H12").Value = 363.6343
$ws.Range("I12").Value = 0.1341789346279459
$ws.Range("J12").Value = 0.1341789346279459
$ws.Range("M12").Value = 29.04767233333333
$ws.Range("N12").Value = 87.143017
$ws.Range("O12").Value = 0.09875811426384234
$ws.Range("P12").Value = 0.09875811426384236
$ws.Range("Q12").Value = 3520.909998520345
$ws.Range("R12").Value = 31688.1899866831
$ws.Range("S12").Value = 0.01325125855778731
$ws.Range("T12").Value = 0.01325125855778731
$ws.Range("G13").Value = 121.2114333333333
$ws.Range("H13").Value = 363.6343
$ws.Range("I13").Value = 0.1341789346279459
$ws.Range("J13").Value = 0.1341789346279459
$ws.Range("M13").Value = 81.87450533333333
$ws.Range("N13").Value = 245.623516
$ws.Range("O13").Value = 0.2783621234850603
$ws.Range("P13").Value = 0.2783621234850603
$ws.Range("Q13").Value = 9924.126144910977
$ws.Range("R13").Value = 89317.1353041988
$ws.Range("S13").Value = 0.03735033316999811
$ws.Range("T13").Value = 0.03735033316999811
$ws.Range("G14").Value = 0.4548216666666667
$ws.Range("H14").Value = 1.364465
$ws.Range("I14").Value = 0.0005034796223489374
$ws.Range("J14").Value = 0.0005034796223489374
$ws.Range("M14").Value = 1.918906333333333
$ws.Range("N14").Value = 5.756718999999999
$ws.Range("O14").Value = 0.006524019162508824
$ws.Range("P14").Value = 0.006524019162508824
$ws.Range("Q14").Value = 0.8727601767038888
$ws.Range("R14").Value = 7.854841590334999
$ws.Range("S14").Value = 0.000003284710704137174
$ws.Range("T14").Value = 0.000003284710704137174
$ws.Range("G15").Value = 0.4548216666666667
$ws.Range("H15").Value = 1.364465
$ws.Range("I15").Value = 0.0005034796223489374
$ws.Range("J15").Value = 0.0005034796223489374
$ws.Range("O15").Value = 0.6163557430885885
$ws.Range("P15").Value = 0.6163557430885885
$ws.Range("Q15").Value = 82.45388829354556
$ws.Range("R15").Value = 742.08499464191
$ws.Range("S15").Value = 0.0003103225567628413
$ws.Range("T15").Value = 0.0003103225567628413
$ws.Range("G16").Value = 0.4548216666666667
$ws.Range("H16").Value = 1.364465
$ws.Range("I16").Value = 0.0005034796223489374
$ws.Range("J16").Value = 0.0005034796223489374
$ws.Range("M16").Value = 29.04767233333333
$ws.Range("N16").Value = 87.143017
$ws.Range("O16").Value = 0.09875811426384234
$ws.Range("P16").Value = 0.09875811426384236
$ws.Range("Q16").Value = 13.21151074343389
$ws.Range("R16").Value = 118.903596690905
$ws.Range("S16").Value = 0.00004972269807345255
$ws.Range("T16").Value = 0.00004972269807345256
$ws.Range("G17").Value = 0.4548216666666667
$ws.Range("H17").Value = 1.364465
$ws.Range("I17").Value = 0.0005034796223489374
$ws.Range("J17").Value = 0.0005034796223489374
$ws.Range("M17").Value = 81.87450533333333
$ws.Range("N17").Value = 245.623516
$ws.Range("O17").Value = 0.2783621234850603
$ws.Range("P17").Value = 0.2783621234850603
$ws.Range("Q17").Value = 9924.126144910977
$ws.Range("R17").Value = 89317.1353041988
$ws.Range("S17").Value = 0.03735033316999811
$ws.Range("T17").Value = 0.03735033316999811
